# MBTW1S1.xlsx — "Add files via upload" edit
#
# 1. Trim the unit suffixes off the Fructose/Biomass/Acetate column headers.
# 2. Move the sheet's active selection from F6 to D1:D2 (D1 active).
# 3. Nudge the saved window's horizontal screen position (xWindow 16560 -> 8500).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text edits (shared strings used by B1/C1/D1) ---
$ws.Range("B1").Value = "Fructose"
$ws.Range("C1").Value = "Biomass"
$ws.Range("D1").Value = "Acetate"

# --- 2. Selection change: activeCell D1, sqref D1:D2 ---
$ws.Range("D1:D2").Select()

# --- 3. Window position (xWindow attribute) ---
# Best-effort; some hosts don't persist window geometry back to the file.
try {
    $excel.ActiveWindow.Left = 8500
} catch {
}
